$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.099.69'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.651.49'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.23%  '
$c = $ws.Range('D5')
$s = $c.Style
$c.Value = "'218.33"
$c.Style = $s
$c = $ws.Range('D6')
$s = $c.Style
$c.Value = "'0.5205"
$c.Style = $s
$ws.Range('E6').Value = '  -0.14%  '
$c = $ws.Range('D7')
$s = $c.Style
$c.Value = "'1.004"
$c.Style = $s
$ws.Range('E7').Value = '  -0.26%  '
$c = $ws.Range('D9')
$s = $c.Style
$c.Value = "'0.06331"
$c.Style = $s
$ws.Range('E9').Value = '  +0.48%  '
$c = $ws.Range('D10')
$s = $c.Style
$c.Value = "'20.37"
$c.Style = $s
$ws.Range('E10').Value = '  -0.81%  '
$c = $ws.Range('D11')
$s = $c.Style
$c.Value = "'0.07693"
$c.Style = $s
$ws.Range('E11').Value = '  -1.51%  '
$c = $ws.Range('D12')
$s = $c.Style
$c.Value = "'4.612"
$c.Style = $s
$ws.Range('E12').Value = '  +2.87%  '
$ws.Range('D13').Value = '1.654.79'
$ws.Range('E13').Value = '  +0.17%  '
$ws.Range('D14').Value = '1.879.54'
$ws.Range('E14').Value = '  +0.07%  '
$c = $ws.Range('D15')
$s = $c.Style
$c.Value = "'0.5591"
$c.Style = $s
$ws.Range('E15').Value = '  +0.66%  '
$ws.Range('D16').Value = '0.0₅8149'
$ws.Range('E16').Value = '  +1.64%  '
$c = $ws.Range('D17')
$s = $c.Style
$c.Value = "'65.33"
$c.Style = $s
$ws.Range('E17').Value = '  +0.66%  '
$ws.Range('D18').Value = '26.103.33'
$ws.Range('E18').Value = '  +0.09%  '
$c = $ws.Range('D19')
$s = $c.Style
$c.Value = "'1.004"
$c.Style = $s
$ws.Range('E19').Value = '  -0.22%  '
$c = $ws.Range('D20')
$s = $c.Style
$c.Value = "'4.627"
$c.Style = $s
$ws.Range('E20').Value = '  -0.32%  '
$ws.Range('E21').Value = '  +4.16%  '
$c = $ws.Range('D22')
$s = $c.Style
$c.Value = "'190.57"
$c.Style = $s
$ws.Range('E22').Value = '  -2.13%  '
$c = $ws.Range('D23')
$s = $c.Style
$c.Value = "'5.931"
$c.Style = $s
$ws.Range('E23').Value = '  -0.37%  '
$c = $ws.Range('D25')
$s = $c.Style
$c.Value = "'144.56"
$c.Style = $s
$ws.Range('E25').Value = '  -1.57%  '
$c = $ws.Range('D26')
$s = $c.Style
$c.Value = "'0.1189"
$c.Style = $s
$ws.Range('E26').Value = '  -1.46%  '
$c = $ws.Range('D27')
$s = $c.Style
$c.Value = "'7.212"
$c.Style = $s
$ws.Range('E27').Value = '  +0.38%  '
$c = $ws.Range('D28')
$s = $c.Style
$c.Value = "'15.92"
$c.Style = $s
$ws.Range('E28').Value = '  +0.08%  '
$c = $ws.Range('D29')
$s = $c.Style
$c.Value = "'1.503"
$c.Style = $s
$ws.Range('E29').Value = '  +2.09%  '
$c = $ws.Range('D30')
$s = $c.Style
$c.Value = "'0.05479"
$c.Style = $s
$ws.Range('E30').Value = '  -3.64%  '
$ws.Range('E31').Value = '  +0.22%  '
$c = $ws.Range('D32')
$s = $c.Style
$c.Value = "'3.444"
$c.Style = $s
$ws.Range('E32').Value = '  -1.05%  '
$c = $ws.Range('D33')
$s = $c.Style
$c.Value = "'3.356"
$c.Style = $s
$ws.Range('E33').Value = '  -0.21%  '
$c = $ws.Range('D34')
$s = $c.Style
$c.Value = "'1.559"
$c.Style = $s
$ws.Range('E34').Value = '  -2.13%  '
$c = $ws.Range('D35')
$s = $c.Style
$c.Value = "'0.9480"
$c.Style = $s
$ws.Range('E35').Value = '  -0.24%  '
$c = $ws.Range('D36')
$s = $c.Style
$c.Value = "'2.783"
$c.Style = $s
$ws.Range('E36').Value = '  -0.64%  '
$c = $ws.Range('D37')
$s = $c.Style
$c.Value = "'2.397"
$c.Style = $s
$ws.Range('E37').Value = '  -0.59%  '
$c = $ws.Range('D38')
$s = $c.Style
$c.Value = "'0.5631"
$c.Style = $s
$ws.Range('E38').Value = '  -0.51%  '
$c = $ws.Range('D39')
$s = $c.Style
$c.Value = "'0.01576"
$c.Style = $s
$ws.Range('E39').Value = '  -0.66%  '
$c = $ws.Range('D40')
$s = $c.Style
$c.Value = "'5.860"
$c.Style = $s
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('E41').Value = '  -0.19%  '
$c = $ws.Range('D42')
$s = $c.Style
$c.Value = "'0.8336"
$c.Style = $s
$ws.Range('E42').Value = '  -0.92%  '
$ws.Range('D43').Value = '1.030.18'
$ws.Range('E43').Value = '  -2.57%  '
$c = $ws.Range('D44')
$s = $c.Style
$c.Value = "'100.87"
$c.Style = $s
$ws.Range('E44').Value = '  -3.81%  '
$ws.Range('D45').Value = '1.792.14'
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₈110'
$ws.Range('E46').Value = '  +5.54%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D47')
$s = $c.Style
$c.Value = "'57.61"
$c.Style = $s
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('E49').Value = '  -0.03%  '
$c = $ws.Range('D50')
$s = $c.Style
$c.Value = "'7.998"
$c.Style = $s
$ws.Range('E50').Value = '  +0.84%  '
$c = $ws.Range('D51')
$s = $c.Style
$c.Value = "'0.05170"
$c.Style = $s
$ws.Range('E51').Value = '  -2.78%  '
